# Weekly update: a new price record was inserted at the top of the
# historical series (rows 92..178), pushing every existing record down by
# one row and appending the previously-last record as a brand new row 179.
#
# Columns D (Fecha), J (Volumen), K (Precio minimo), M (Precio promedio
# ponderado) and P (Precio $/Kg) are the only ones that vary per record in
# this block, so those are the only ones that need to shift; every other
# column (A,B,C,E,F,G,H,I,L,N,O,Q,R) is constant across the block and is
# simply copied along with the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 92
$lastRow = 178
$newLastRow = $lastRow + 1

# New values being inserted at the front of the block (row 92).
$newDate   = 44512
$newVolume = 3340
$newMin    = 400
$newAvg    = 450
$newKg     = 900

# --- Capture the "before" state of the shifting columns for every row in
#     the block, plus the full static row template, before anything is
#     overwritten. ---
$vals = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $row = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        M = $ws.Cells.Item($r, 13).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
    $vals[$r] = $row
}

# --- Build the brand new last row (179) by duplicating every static
#     column from the old last row (178), then stamp in the shifted
#     values (which equal the old row 178's values). ---
for ($c = 1; $c -le 18; $c++) {
    $src = $ws.Cells.Item($lastRow, $c)
    $dst = $ws.Cells.Item($newLastRow, $c)
    $dst.Value = $src.Value2
}
# Column D (Fecha) carries a date number format in this block; replicate it
# on the freshly created row so it matches its neighbours.
$ws.Cells.Item($newLastRow, 4).NumberFormat = $ws.Cells.Item($lastRow, 4).NumberFormat
$ws.Cells.Item($newLastRow, 4).Value  = $vals[$lastRow].D
$ws.Cells.Item($newLastRow, 10).Value = $vals[$lastRow].J
$ws.Cells.Item($newLastRow, 11).Value = $vals[$lastRow].K
$ws.Cells.Item($newLastRow, 13).Value = $vals[$lastRow].M
$ws.Cells.Item($newLastRow, 16).Value = $vals[$lastRow].P

# --- Shift rows 178 down to 93: each row takes on the value that used to
#     belong to the row right above it. ---
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $prev = $vals[$r - 1]
    $ws.Cells.Item($r, 4).Value  = $prev.D
    $ws.Cells.Item($r, 10).Value = $prev.J
    $ws.Cells.Item($r, 11).Value = $prev.K
    $ws.Cells.Item($r, 13).Value = $prev.M
    $ws.Cells.Item($r, 16).Value = $prev.P
}

# --- Row 92 gets the brand new record. ---
$ws.Cells.Item($firstRow, 4).Value  = $newDate
$ws.Cells.Item($firstRow, 10).Value = $newVolume
$ws.Cells.Item($firstRow, 11).Value = $newMin
$ws.Cells.Item($firstRow, 13).Value = $newAvg
$ws.Cells.Item($firstRow, 16).Value = $newKg
